$d = $word.ActiveDocument
for ($i=1; $i -le 18; $i++) {
    try {
        $sr = $d.StoryRanges.Item($i)
        Write-Output ($i.ToString() + ": len=" + $sr.Text.Length + " text=[" + $sr.Text.Substring(0, [Math]::Min(40, $sr.Text.Length)) + "]")
    } catch {
        Write-Output ($i.ToString() + ": ERROR " + $_)
    }
}
